$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B3:E5")
$rng.NumberFormat = "@"

# Row 3 -> Ochieng Charles
$ws.Range("A3").Value = "Ochieng Charles"
$ws.Range("B3").Value = "3.00"
$ws.Range("C3").Value = "25.00"
$ws.Range("D3").Value = "-22.00"
$ws.Range("E3").Value = "12.00%"

# Row 4 -> Lenah Cheloti
$ws.Range("A4").Value = "Lenah Cheloti"
$ws.Range("B4").Value = "5.00"
$ws.Range("C4").Value = "18.00"
$ws.Range("D4").Value = "-13.00"
$ws.Range("E4").Value = "27.78%"

# Row 5 -> Moses  Ngugi
$ws.Range("A5").Value = "Moses  Ngugi"
$ws.Range("B5").Value = "4.00"
$ws.Range("C5").Value = "30.00"
$ws.Range("D5").Value = "-26.00"
$ws.Range("E5").Value = "13.33%"
